$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.225
$ws.Range("C2").Value = 0.425
$ws.Range("P2").Value = 0.25
$ws.Range("S2").Value = 0.1
$ws.Range("P3").Value = 0.6666666666666666
$ws.Range("S3").Value = 0.3333333333333333
$ws.Range("J4").Value = 0.1
$ws.Range("P4").Value = 0.8
$ws.Range("S4").Value = 0.1
$ws.Range("F6").Value = 0.1
$ws.Range("J6").Value = 0.2
$ws.Range("Q6").Value = 0.25
$ws.Range("R6").Value = 0.1
$ws.Range("S6").Value = 0.35
$ws.Range("Q7").Value = 0.125
$ws.Range("S7").Value = 0.125
$ws.Range("B8").Value = 0.07142857142857142
$ws.Range("D8").Value = 0.1071428571428571
$ws.Range("F8").Value = 0.07142857142857142
$ws.Range("J8").Value = 0.1071428571428571
$ws.Range("O8").Value = 0.03571428571428571
$ws.Range("Q8").Value = 0.2142857142857143
$ws.Range("R8").Value = 0.1785714285714286
$ws.Range("S8").Value = 0.2142857142857143
$ws.Range("B9").Value = 0.2222222222222222
$ws.Range("F9").Value = 0.07407407407407407
$ws.Range("J9").Value = 0.03703703703703703
$ws.Range("Q9").Value = 0.2962962962962963
$ws.Range("R9").Value = 0.1481481481481481
$ws.Range("S9").Value = 0.2222222222222222
$ws.Range("B10").Value = 0.1020408163265306
$ws.Range("D10").Value = 0.04591836734693878
$ws.Range("F10").Value = 0.04591836734693878
$ws.Range("J10").Value = 0.1224489795918367
$ws.Range("O10").Value = 0.01530612244897959
$ws.Range("Q10").Value = 0.3061224489795918
$ws.Range("R10").Value = 0.1020408163265306
$ws.Range("S10").Value = 0.2602040816326531
$ws.Range("J11").Value = 0.2
$ws.Range("K11").Value = 0.1333333333333333
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.6
$ws.Range("J12").Value = 0.4
$ws.Range("G13").Value = 0.75
$ws.Range("S13").Value = 0.25
$ws.Range("F15").Value = 0.06896551724137931
$ws.Range("H15").Value = 0.103448275862069
$ws.Range("I15").Value = 0.03448275862068965
$ws.Range("J15").Value = 0.5517241379310345
$ws.Range("M15").Value = 0.03448275862068965
$ws.Range("O15").Value = 0.03448275862068965
$ws.Range("S15").Value = 0.1724137931034483
$ws.Range("H16").Value = 0.1
$ws.Range("I16").Value = 0.1
$ws.Range("J16").Value = 0.6333333333333333
$ws.Range("K16").Value = 0.03333333333333333
$ws.Range("O16").Value = 0.03333333333333333
$ws.Range("S16").Value = 0.1
$ws.Range("H17").Value = 0.1052631578947368
$ws.Range("I17").Value = 0.07894736842105263
$ws.Range("J17").Value = 0.5526315789473685
$ws.Range("K17").Value = 0.07894736842105263
$ws.Range("O17").Value = 0.09210526315789473
$ws.Range("S17").Value = 0.09210526315789473
$ws.Range("H18").Value = 0.06060606060606061
$ws.Range("I18").Value = 0.06060606060606061
$ws.Range("J18").Value = 0.5757575757575758
$ws.Range("M18").Value = 0.0303030303030303
$ws.Range("O18").Value = 0.1212121212121212
$ws.Range("S18").Value = 0.1515151515151515
$ws.Range("F19").Value = 0.01834862385321101
$ws.Range("H19").Value = 0.1192660550458716
$ws.Range("I19").Value = 0.1376146788990826
$ws.Range("J19").Value = 0.5229357798165137
$ws.Range("K19").Value = 0.03669724770642202
$ws.Range("M19").Value = 0.01834862385321101
$ws.Range("O19").Value = 0.07339449541284404
$ws.Range("S19").Value = 0.07339449541284404
